$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2=3, B2=262 (was A2=0, B2=409)
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 262

# New row 3: A3=0, B3=208
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 208

# New row 4 (was old row3 shifted down): A4=1, B4=131
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 131

# New row 5: A5=2, B5=113
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 113

# Copy the style of A2 (already style s="1") onto the new cells A3:A5
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
